# Generate Report for Handback
# Refresh the "Correspond Handoff Datetime" (D2/D3) and
# "Correspond Handback DateTime" (G2/G3) timestamps on the per-language
# report sheets to reflect the latest handback run.

$wb = $excel.ActiveWorkbook

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("D2").Value = "2016-02-23 08:17:24"
$wsZh.Range("D3").Value = "2016-02-23 08:17:24"
$wsZh.Range("G2").Value = "2016-02-23 08:18:29"
$wsZh.Range("G3").Value = "2016-02-23 08:18:29"

$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("D2").Value = "2016-02-23 08:17:40"
$wsDe.Range("D3").Value = "2016-02-23 08:17:40"
$wsDe.Range("G2").Value = "2016-02-23 08:18:56"
$wsDe.Range("G3").Value = "2016-02-23 08:18:56"
